$d = $word.ActiveDocument

# Locate the "LOB1018: Física I (Requisito)" paragraph (end of the Requisitos
# section). Immediately after it, the document currently has three trailing
# paragraphs that must be removed:
#   1. an empty paragraph
#   2. "Ver no Jupiter Salvar em pdf Salvar em docx"
#   3. "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github
#      pages. Original theme under Creative Commons Attribution"
# The empty paragraph and the final page-break paragraph that follow must be
# kept untouched.
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*LOB1018: Física I (Requisito)*") {
        $anchorIndex = $i
        break
    }
}

if ($anchorIndex -gt 0) {
    $firstToRemove = $d.Paragraphs.Item($anchorIndex + 1)
    $lastToRemove = $d.Paragraphs.Item($anchorIndex + 3)

    $deleteRange = $d.Range($firstToRemove.Range.Start, $lastToRemove.Range.End)
    $deleteRange.Delete()
}
